$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 325.13333
$ws.Range("I39").Value = 79.09090999999999
$ws.Range("J39").Value = 1001.75
$ws.Range("K39").Value = 237.27273
$ws.Range("L39").Value = 3005.25
$ws.Range("M39").Value = 58.72727000000003
$ws.Range("N39").Value = -3597.25
$ws.Range("H40").Value = 876
$ws.Range("I40").Value = 716.9231
$ws.Range("J40").Value = 1171.4286
$ws.Range("K40").Value = 716.9231
$ws.Range("L40").Value = 1171.4286
$ws.Range("M40").Value = -541.9231
$ws.Range("N40").Value = -1521.4286
$ws.Range("H52").Value = 8000
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 8000
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 24000
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = -24320
$ws.Range("H58").Value = 671836
$ws.Range("J58").Value = 916012.75
$ws.Range("L58").Value = 2748038.25
$ws.Range("N58").Value = -2748338.25
$ws.Range("H64").Value = 3700
$ws.Range("I64").Value = 3320
$ws.Range("J64").Value = 4016.6667
$ws.Range("K64").Value = 3320
$ws.Range("L64").Value = 4016.6667
$ws.Range("M64").Value = -3072
$ws.Range("N64").Value = -4512.6667
$ws.Range("H67").Value = 3700
$ws.Range("I67").Value = 3320
$ws.Range("J67").Value = 4016.6667
$ws.Range("K67").Value = 3320
$ws.Range("L67").Value = 4016.6667
$ws.Range("M67").Value = -2462
$ws.Range("N67").Value = -5732.6667
$ws.Range("H129").Value = 418830.6
$ws.Range("J129").Value = 478625.56
$ws.Range("L129").Value = 1435876.68
$ws.Range("N129").Value = -1445876.68
$ws.Range("H138").Value = 3398.0908
$ws.Range("I138").Value = 2132.3684
$ws.Range("J138").Value = 4360.04
$ws.Range("K138").Value = 6397.1052
$ws.Range("L138").Value = 13080.12
$ws.Range("M138").Value = -1257.1052
$ws.Range("N138").Value = -23360.12
$ws.Range("H141").Value = 1568.3778
$ws.Range("I141").Value = 1380.4048
$ws.Range("J141").Value = 4200
$ws.Range("K141").Value = 4141.2144
$ws.Range("L141").Value = 12600
$ws.Range("M141").Value = 1038.7856
$ws.Range("N141").Value = -22960

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1101.3077
$ws.Range("I110").Value = 1101.3077
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 1101.3077
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 943.6922999999999
$ws.Range("N110").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2744.5625
$ws.Range("J86").Value = 3835.6667
$ws.Range("L86").Value = 3835.6667
$ws.Range("N86").Value = -6081.6667
$ws.Range("H89").Value = 2744.5625
$ws.Range("J89").Value = 3835.6667
$ws.Range("L89").Value = 19178.3335
$ws.Range("N89").Value = -30410.3335
$ws.Range("H94").Value = 907.6
$ws.Range("I94").Value = 511.3684
$ws.Range("J94").Value = 1592
$ws.Range("K94").Value = 511.3684
$ws.Range("L94").Value = 1592
$ws.Range("M94").Value = -60.36840000000001
$ws.Range("N94").Value = -2494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 180.3125
$ws.Range("J22").Value = 250
$ws.Range("L22").Value = 250
$ws.Range("N22").Value = -950
$ws.Range("H31").Value = 4283.4526
$ws.Range("J31").Value = 5902.1377
$ws.Range("L31").Value = 5902.1377
$ws.Range("N31").Value = -6492.1377
$ws.Range("H34").Value = 4283.4526
$ws.Range("J34").Value = 5902.1377
$ws.Range("L34").Value = 5902.1377
$ws.Range("N34").Value = -6306.1377
$ws.Range("H132").Value = 31252058
$ws.Range("I132").Value = 34484092
$ws.Range("K132").Value = 103452276
$ws.Range("M132").Value = -103449746

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 115
$ws.Range("I4").Value = 119.42857
$ws.Range("J4").Value = 99.5
$ws.Range("K4").Value = 358.28571
$ws.Range("L4").Value = 298.5
$ws.Range("M4").Value = -246.28571
$ws.Range("N4").Value = -522.5
$ws.Range("H22").Value = 20640
$ws.Range("J22").Value = 800
$ws.Range("L22").Value = 2400
$ws.Range("N22").Value = -2738
$ws.Range("H27").Value = 20640
$ws.Range("J27").Value = 800
$ws.Range("L27").Value = 2400
$ws.Range("N27").Value = -2604
$ws.Range("H122").Value = 1013.1429
$ws.Range("J122").Value = 1013.1429
$ws.Range("L122").Value = 9118.286100000001
$ws.Range("N122").Value = -14018.2861
$ws.Range("H131").Value = 697.9400000000001
$ws.Range("I131").Value = 380.69232
$ws.Range("J131").Value = 745.34485
$ws.Range("K131").Value = 1142.07696
$ws.Range("L131").Value = 2236.03455
$ws.Range("M131").Value = 3897.92304
$ws.Range("N131").Value = -12316.03455
$ws.Range("H136").Value = 3062.5186
$ws.Range("J136").Value = 4991.385
$ws.Range("L136").Value = 14974.155
$ws.Range("N136").Value = -25174.155

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 981
$ws.Range("I97").Value = 227.25
$ws.Range("K97").Value = 227.25
$ws.Range("M97").Value = 268.75
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("H126").Value = 3460.111
$ws.Range("I126").Value = 2428.2778
$ws.Range("J126").Value = 5523.778
$ws.Range("K126").Value = 7284.8334
$ws.Range("L126").Value = 16571.334
$ws.Range("M126").Value = -4814.8334
$ws.Range("N126").Value = -21511.334
$ws.Range("H127").Value = 20000
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 20000
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 20000
$ws.Range("N127").Value = -29920
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("H129").Value = 49999
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 49999
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 49999
$ws.Range("N129").Value = -59999
$ws.Range("H130").Value = 48153.26
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 48153.26
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 48153.26
$ws.Range("N130").Value = -58193.26
$ws.Range("H131").Value = 40347.383
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 40347.383
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 40347.383
$ws.Range("N131").Value = -50427.383
$ws.Range("H132").Value = 4401073.5
$ws.Range("I132").Value = 6688954.5
$ws.Range("J132").Value = 54099.6
$ws.Range("K132").Value = 20066863.5
$ws.Range("L132").Value = 162298.8
$ws.Range("M132").Value = -20064333.5
$ws.Range("N132").Value = -167358.8
$ws.Range("H133").Value = 40250
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 40250
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 40250
$ws.Range("N133").Value = -50370
$ws.Range("H134").Value = 20265.2
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 20265.2
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 60795.60000000001
$ws.Range("N134").Value = -65865.60000000001
$ws.Range("H135").Value = 37207.8
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 37207.8
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 37207.8
$ws.Range("N135").Value = -47347.8
$ws.Range("H136").Value = 16600.2
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 16600.2
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 49800.60000000001
$ws.Range("N136").Value = -54900.60000000001
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("H139").Value = 39413.5
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 39413.5
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 39413.5
$ws.Range("N139").Value = -49693.5
$ws.Range("H140").Value = 60700
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 60700
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 60700
$ws.Range("N140").Value = -71060
$ws.Range("H141").Value = 67822.8
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 67822.8
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 67822.8
$ws.Range("N141").Value = -78182.8
